$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at row 388 (shifts existing rows 388..480 down to 392..484)
$ws.Range("A388:A391").EntireRow.Insert()

# Fill in the 4 new rows with the new weekly price data (date 44543 = 2021-12-13)
$newRows = @(
    @{ Row=388; D=44543; I="1a (cosecha)"; J=600;  K=14000; L=14000; M=14000; O="Provincia de Talagante"; P=1400 },
    @{ Row=389; D=44543; I="2a (cosecha)"; J=400;  K=13000; L=13000; M=13000; O="Provincia de Talagante"; P=1300 },
    @{ Row=390; D=44543; I="3a (cosecha)"; J=200;  K=11000; L=11000; M=11000; O="Provincia de Talagante"; P=1100 },
    @{ Row=391; D=44543; I="Primera";      J=1800; K=16500; L=17000; M=16639; O="China";                  P=1664 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 6
    $ws.Cells.Item($row, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 100112003
    $ws.Cells.Item($row, 7).Value = "Ajo"
    $ws.Cells.Item($row, 8).Value = "Chino"
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "$/caja 10 kilos"
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 10
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
